$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$row = 80

$ws.Cells.Item($row, 1).Value = 1
$ws.Cells.Item($row, 2).Value = "Agrícola del Norte S.A. de Arica"
$ws.Cells.Item($row, 3).Value = "Arica y Parinacota"
$ws.Cells.Item($row, 4).Value = 44832
$ws.Cells.Item($row, 4).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item($row, 5).Value = 15
$ws.Cells.Item($row, 6).Value = 100112012
$ws.Cells.Item($row, 7).Value = "Espinaca"
$ws.Cells.Item($row, 8).Value = "Sin especificar"
$ws.Cells.Item($row, 9).Value = "Primera"
$ws.Cells.Item($row, 10).Value = 300
$ws.Cells.Item($row, 11).Value = 2000
$ws.Cells.Item($row, 12).Value = 2500
$ws.Cells.Item($row, 13).Value = 2250
$ws.Cells.Item($row, 14).Value = "$/atado 2,5 a 3 kilos"
$ws.Cells.Item($row, 15).Value = "Región de Arica y Parinacota"
$ws.Cells.Item($row, 16).Value = 750
$ws.Cells.Item($row, 17).Value = 3
$ws.Cells.Item($row, 18).Value = "Hortaliza"
